$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force affected cells to remain plain text so numeric-looking values
# (e.g. "50.066.76", "0.530", "5.11") are not reinterpreted/rounded as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '50.066.76'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.669.37'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +7.53%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '114.38'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +8.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '327.41'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.530'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.28%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.559'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.37'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.18'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0827'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.53%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.081.37'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +7.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.687.29'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +8.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.879'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +6.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '50.031.41'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +4.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.33'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +5.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.82'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.94'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0962'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.64%  '
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '285.82'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.80%  '
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.88'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.60'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.01'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +5.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.10%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.89'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +7.33%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.24'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.80%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.41'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.19%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +5.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.79'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0818'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.11'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +13.33%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +8.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.14'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +9.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.39'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.53%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.62%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.20'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0320'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +6.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.112.71'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +5.86%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +6.39%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +14.14%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +4.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.14'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.39'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.53'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +6.16%  '
